# "seteo fecha actual pago" -- set the current payment date, clarify the
# payment label/value on the IMPORTARHOJA sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IMPORTARHOJA")

# Clarify the "pago" label
$ws.Range("A2").Value = "Pagó (Si, No):"

# Record that it was paid ("si") instead of the stray numeric leftover
$ws.Range("B2").Value = "si"

# Set the payment date (D1) to the current date
$ws.Range("D1").Value = Get-Date -Year 2014 -Month 12 -Day 19 -Hour 0 -Minute 0 -Second 0

# Move the active selection to B3, matching where the user left off
$ws.Range("B3").Select()
